# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.225.60"
$ws.Range("E2").Value = "  +6.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.680.95"
$ws.Range("E3").Value = "  +19.12%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.91"
$ws.Range("E5").Value = "  +3.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.03"
$ws.Range("E6").Value = "  +6.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.679.88"
$ws.Range("E7").Value = "  +19.15%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +4.18%  "

$ws.Range("E10").Value = "  +8.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  +3.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").Value = "  +6.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.88"
$ws.Range("E13").Value = "  +11.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  +6.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.301.49"
$ws.Range("E15").Value = "  +19.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.684.54"
$ws.Range("E16").Value = "  +19.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.257.88"
$ws.Range("E17").Value = "  +7.04%  "

$ws.Range("E18").Value = "  +1.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  +7.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.87"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "515.43"
$ws.Range("E21").Value = "  +6.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.17"
$ws.Range("E22").Value = "  +18.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.745"
$ws.Range("E23").Value = "  +8.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.50"
$ws.Range("E24").Value = "  +5.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  +8.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.40"
$ws.Range("E26").Value = "  +5.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("E27").Value = "  +8.09%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +12.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.13"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.66"
$ws.Range("E31").Value = "  +13.35%  "

$ws.Range("E32").Value = "  +7.22%  "

$ws.Range("E33").Value = "  +17.68%  "

$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +9.80%  "

$ws.Range("E37").Value = "  +7.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.342"
$ws.Range("E38").Value = "  +10.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.15"
$ws.Range("E39").Value = "  +9.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.87"
$ws.Range("E40").Value = "  +3.74%  "

$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.24"
$ws.Range("E42").Value = "  -6.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.168.71"
$ws.Range("E43").Value = "  +14.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.80"
$ws.Range("E44").Value = "  +6.39%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "407.52"
$ws.Range("E45").Value = "  +10.66%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  +6.11%  "

$ws.Range("E47").Value = "  +6.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.16"
$ws.Range("E48").Value = "  +15.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.61"
$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("E51").Value = "  +12.57%  "
